# Week 17 data logging + Simulate_Season.py tiebreak fix
# Appends this week's per-drive/per-play yardage samples to the running
# logs on the YDS sheet, appends this week's special-teams field-position
# samples on the ST sheet, and rolls the week's totals into the OFF/DEF
# summary sheets and the TURNS sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append new samples to the four running log strings
# ---------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value2 = $wsYDS.Range("B2").Value2 + " 4 9 -1 2 14 3 2 6 4 3 1 6 9 6 0 2 0 -1 1 2 -4 2 7 4 5 -4 2 5 7 2 4 1 3"
$wsYDS.Range("C2").Value2 = $wsYDS.Range("C2").Value2 + " 3 6 7 -3 5 5 4 5 2 0 1 0 6 12 4 7 4 -4 -4 2 2 0 4 5 2 5"
$wsYDS.Range("B3").Value2 = $wsYDS.Range("B3").Value2 + " 10 8 10 11 8 3 6 8 3 5 9 3 8 7 40 16 5 5 45 8 8"
$wsYDS.Range("C3").Value2 = $wsYDS.Range("C3").Value2 + " 15 0 44 18 -1 24 8 35 7 9 5 11 6 2 16 7 34 5"

# ---------------------------------------------------------------------
# OFF sheet: roll Week 17 totals into Home (row 2) / Road (row 3)
# ---------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value = 412
$wsOFF.Range("D2").Value = 28
$wsOFF.Range("F2").Value = 106
$wsOFF.Range("G2").Value = 117
$wsOFF.Range("J2").Value = 56
$wsOFF.Range("L2").Value = 678
$wsOFF.Range("M2").Value = 450
$wsOFF.Range("O2").Value = 51
$wsOFF.Range("P2").Value = 29
$wsOFF.Range("Q2").Value = 1166

$wsOFF.Range("C3").Value = 373
$wsOFF.Range("E3").Value = 74
$wsOFF.Range("F3").Value = 258
$wsOFF.Range("G3").Value = 76
$wsOFF.Range("I3").Value = 121
$wsOFF.Range("J3").Value = 129

# ---------------------------------------------------------------------
# DEF sheet: roll Week 17 totals into Home (row 2) / Road (row 3)
# ---------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("B2").Value = 7
$wsDEF.Range("C2").Value = 413
$wsDEF.Range("E2").Value = 26
$wsDEF.Range("F2").Value = 121
$wsDEF.Range("G2").Value = 130
$wsDEF.Range("I2").Value = 17
$wsDEF.Range("J2").Value = 72
$wsDEF.Range("L2").Value = 576
$wsDEF.Range("M2").Value = 354
$wsDEF.Range("O2").Value = 36
$wsDEF.Range("P2").Value = 16
$wsDEF.Range("Q2").Value = 1061

$wsDEF.Range("C3").Value = 331
$wsDEF.Range("D3").Value = 9
$wsDEF.Range("E3").Value = 56
$wsDEF.Range("F3").Value = 190
$wsDEF.Range("G3").Value = 78
$wsDEF.Range("I3").Value = 109
$wsDEF.Range("J3").Value = 104
$wsDEF.Range("N3").Value = 33

# ---------------------------------------------------------------------
# ST sheet: update KO/PT summary counts and append new PT/KO samples
# ---------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 174
$wsST.Range("D2").Value = 103
$wsST.Range("F2").Value = 222
$wsST.Range("G2").Value = 209
$wsST.Range("J2").Value = 109
$wsST.Range("K2").Value = 104
$wsST.Range("B3").Value = 109

$wsST.Range("B6").Value2 = $wsST.Range("B6").Value2 + " 47 1"
$wsST.Range("D3").Value2 = $wsST.Range("D3").Value2 + " 25 48"
$wsST.Range("D4").Value2 = $wsST.Range("D4").Value2 + " 0 4"
$wsST.Range("D5").Value2 = $wsST.Range("D5").Value2 + " 0 1 4"

# ---------------------------------------------------------------------
# TURNS sheet: fixed tiebreaking-related giveaway count
# ---------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("D3").Value = 11
